$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'260.55"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'1.47%"
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'27.22"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'2.36%"
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'4.685"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'0.78%"
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'0.06153"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'3.77%"
$ws.Range("E5").ClearFormats()
$ws.Range("E6").Value = "'0.80%"
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'0.8511"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'-0.55%"
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.9154"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'0.42%"
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.1406"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'2.10%"
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.04722"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'9.94%"
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.07088"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'1.12%"
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.03101"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'2.59%"
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'0.09054"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'-0.58%"
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'0.001539"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'0.76%"
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'0.0006149"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'1.86%"
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'0.006058"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'0.10%"
$ws.Range("E16").ClearFormats()
$ws.Range("D18").Value = "'3.158"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'0.70%"
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = "'1.38%"
$ws.Range("E19").ClearFormats()
$ws.Range("E20").Value = "'0.94%"
$ws.Range("E20").ClearFormats()
$ws.Range("E21").Value = "'1.34%"
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'4.080"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'5.08%"
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'0.04225"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'0.48%"
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'0.001215"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").Value = "'0.003800"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'-18.31%"
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "'0.21%"
$ws.Range("E26").ClearFormats()
$ws.Range("E27").Value = "'-8.13%"
$ws.Range("E27").ClearFormats()
$ws.Range("D40").Value = "'0.03873"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'1.87%"
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.1112"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'1.08%"
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'0.004074"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'10.01%"
$ws.Range("E42").ClearFormats()
$ws.Range("E43").Value = "'13.84%"
$ws.Range("E43").ClearFormats()
$ws.Range("E44").Value = "'-9.81%"
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.00005165"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'0.56%"
$ws.Range("E45").ClearFormats()
$ws.Range("E46").Value = "'0.06%"
$ws.Range("E46").ClearFormats()
$ws.Range("E47").Value = "'8.06%"
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'0.1624"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'-32.59%"
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'0.06%"
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'0.06%"
$ws.Range("E50").ClearFormats()
